$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "16/3/2010" attendance column (column L) for rows 3-8
$ws.Range("L3").Value = "-"
$ws.Range("L4").Value = "x"
$ws.Range("L5").Value = "x"
$ws.Range("L6").Value = "-"
$ws.Range("L7").Value = "x"
$ws.Range("L8").Value = "x"

# Update the selected cell to L9
$ws.Range("L9").Select()
